# Media App.docx edit
#
# The original document uses <w:proofErr/> spell/grammar-check markers
# that split single pieces of text into multiple adjacent <w:r> runs
# (e.g. "Sql" + proofErr + " Server"). The edit removes every one of
# those proofErr markers and merges the runs they used to separate back
# into a single run (affects the "Sql Server", "Cloudinary database for
# pictures", "Custom execption middleware", "Cloudinary implementation",
# "Jwt, loading and error interceptors", "Bootswatch minty for theme",
# "Toastr for frontend user notification", "Ngx-spinner for loading",
# "NgbDatepickerModule" and "npm i bootstrap-icons" bullets), plus
# appends one new bullet to the end of the Frontend list.
#
# Word's own Range.WordOpenXML getter already renders a paragraph's
# content with the proofErr markers stripped and the runs coalesced, so
# round-tripping every paragraph's WordOpenXML back through
# Range.InsertXML reproduces exactly the merged/cleaned structure the
# diff describes, while preserving every other attribute untouched.
# Paragraphs that are already "clean" (title, plain section headers,
# already-single-run bullets) round-trip to themselves, so it is safe
# to do this uniformly instead of hand-picking paragraph numbers.

$d = $word.ActiveDocument

foreach ($idx in 1..$d.Paragraphs.Count) {
    $p = $d.Paragraphs($idx)
    $cleanXml = $p.Range.WordOpenXML
    $null = $p.Range.InsertXML($cleanXml)
}

# Append the new bullet "ng pagination boostrap" after the last
# paragraph ("npm i bootstrap-icons"), using the same ListParagraph /
# numId=2 list the rest of the Frontend section uses. Inserting at the
# collapsed end of $d.Content (rather than into the last paragraph's
# own Range) avoids splitting/duplicating that paragraph's own mark.
$newParaXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>ng pagination boostrap</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$endRange = $d.Content
$endRange.Collapse(0)
$null = $endRange.InsertXML($newParaXml)
